# Add the new "Pid4CatRecord" worksheet at the end of the workbook and
# populate its header row, mirroring the auto-generated model sheets
# already present in the workbook.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() (no args) inserts before the active sheet, i.e. at the
# front of the tab order. We need the new sheet appended at the very end,
# so explicitly add it After the current last worksheet.
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Pid4CatRecord"

$headers = @(
    "landing_page_url",
    "status",
    "schema_version",
    "metadata_license",
    "curation_contact",
    "resource_info",
    "related_identifiers",
    "change_log"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
